$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "CMS(e) " label to "CMS(e)(2010)" for every cell that shares it
# (L2:L7 all reference the same shared string), so the shared-string table
# keeps a single consolidated entry with the new text instead of forking a
# new one.
$ws.Range("L2:L7").Value = "CMS(e)(2010)"

# Update the active selection to N6 (previously L8).
$ws.Range("N6").Select()
